# Se agrega que devuelva toda la clase
$wb = $excel.ActiveWorkbook

# Work on the "Aerobox" sheet
$ws = $wb.Worksheets.Item("Aerobox")

# Fill column D (rows 1-3) with "SUAREZ JULIETA" so the class returns full capacity
$ws.Range("D1:D3").ClearFormats()
$ws.Range("D1").Value = "SUAREZ JULIETA"
$ws.Range("D2").Value = "SUAREZ JULIETA"
$ws.Range("D3").Value = "SUAREZ JULIETA"

# Update the active selection on each sheet to match the new state
$ws.Range("D3").Select()

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("N2").Select()

# Leave Aerobox as the active sheet (it was already the selected tab)
$ws.Activate()
